$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'

$ws.Range('D2').Value = '26.595.40'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.594.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '210.85'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.0615'
$ws.Range('D9').Value = '0.245'
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('D11').Value = '0.0836'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '1.817.37'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '1.600.55'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '64.43'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '26.575.53'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '207.79'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = '6.91'
$ws.Range('E21').Value = '  +2.55%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').Value = '145.13'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').Value = '15.21'
$ws.Range('D30').Value = '0.0504'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = '0.654'
$ws.Range('E33').Value = '  -1.37%  '
$ws.Range('D34').Value = '2.92'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').Value = '1.279.94'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('D39').Value = '0.839'
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').Value = '5.46'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').Value = '0.783'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').Value = '63.82'
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('E45').Value = '  +9.52%  '
$ws.Range('D46').Value = '1.730.00'
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').Value = '89.39'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').Value = '1.59'
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('E50').Value = '  +3.88%  '
$ws.Range('E51').Value = '  +0.66%  '
